$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 2635
$ws.Range("F7").Value = 186
$ws.Range("F10").Value = 5747
$ws.Range("F11").Value = 897
$ws.Range("F12").Value = 1502
$ws.Range("F13").Value = 1433
$ws.Range("F14").Value = 617
$ws.Range("F15").Value = 7015
$ws.Range("F17").Value = 58
$ws.Range("F18").Value = 68
$ws.Range("F19").Value = 4782
$ws.Range("F21").Value = 81
$ws.Range("F22").Value = 2409
$ws.Range("F23").Value = 1286
$ws.Range("F24").Value = 466
$ws.Range("F25").Value = 1178
$ws.Range("F26").Value = 246
$ws.Range("F27").Value = 98
$ws.Range("F28").Value = 96
$ws.Range("F29").Value = 181
$ws.Range("F30").Value = 375
$ws.Range("F31").Value = 1311
$ws.Range("F32").Value = 2019
$ws.Range("F33").Value = 257
$ws.Range("F35").Value = 29
$ws.Range("F36").Value = 217
$ws.Range("F37").Value = 1402
$ws.Range("F39").Value = 99
$ws.Range("F40").Value = 534
$ws.Range("F41").Value = 200
$ws.Range("F42").Value = 1661
$ws.Range("F43").Value = 2449
$ws.Range("F45").Value = 85
$ws.Range("F46").Value = 241
$ws.Range("F47").Value = 85
$ws.Range("F48").Value = 37
$ws.Range("F49").Value = 71

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 428
$ws.Range("F5").Value = 459
$ws.Range("F12").Value = 282
$ws.Range("F14").Value = 58
$ws.Range("F15").Value = 192
$ws.Range("F19").Value = 141
$ws.Range("F20").Value = 36
$ws.Range("F26").Value = 310
$ws.Range("F27").Value = 301
$ws.Range("F34").Value = 17
$ws.Range("F36").Value = 1
$ws.Range("F37").Value = 4

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 1674
$ws.Range("F7").Value = 547
$ws.Range("F8").Value = 1365
$ws.Range("F9").Value = 1202
$ws.Range("F10").Value = 1771
$ws.Range("F11").Value = 2286
$ws.Range("F12").Value = 725
$ws.Range("F13").Value = 597

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 1674
$ws.Range("F6").Value = 547
$ws.Range("F7").Value = 2635
$ws.Range("F8").Value = 186
$ws.Range("F9").Value = 1365
$ws.Range("F11").Value = 2286
$ws.Range("F12").Value = 5747
$ws.Range("F13").Value = 725
$ws.Range("F14").Value = 1502
$ws.Range("F15").Value = 1433
$ws.Range("F17").Value = 597
$ws.Range("F19").Value = 58
$ws.Range("F20").Value = 4782
$ws.Range("F21").Value = 2409
$ws.Range("F22").Value = 1286
$ws.Range("F23").Value = 466
$ws.Range("F24").Value = 1178
$ws.Range("F25").Value = 246
$ws.Range("F26").Value = 96
$ws.Range("F27").Value = 282
$ws.Range("F28").Value = 181
$ws.Range("F29").Value = 58
$ws.Range("F30").Value = 192
$ws.Range("F31").Value = 375
$ws.Range("F32").Value = 2019
$ws.Range("F33").Value = 257
$ws.Range("F35").Value = 141
$ws.Range("F36").Value = 36
$ws.Range("F37").Value = 1402
$ws.Range("F39").Value = 534
$ws.Range("F41").Value = 310
$ws.Range("F42").Value = 200
$ws.Range("F44").Value = 1661
$ws.Range("F45").Value = 2449
$ws.Range("F46").Value = 85
$ws.Range("F47").Value = 241
$ws.Range("F48").Value = 85
$ws.Range("F49").Value = 37

